$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (edit in-place within the rich-text shared strings) ---
# "Volume 30   Number  35" -> "Volume 30   Number  36"
$ws.Range("A8").Characters(21, 2).Text = "36"

# "Report Covering the Week  8/28/2023  Through  9/3/2023"
# -> "Report Covering the Week  9/4/2023  Through  9/10/2023"
$ws.Range("C9").Characters(27, 9).Text = "9/4/2023"
$ws.Range("C9").Characters(46, 8).Text = "9/10/2023"

# --- Rows 15 and 26 (Rape / UCR Rape*): the Week-to-Date columns C/D/E change
# shape -- C goes from a number to the "0" placeholder text, D goes from the
# "0" placeholder text to a number, and E goes from the "***.*" placeholder
# text to a computed percentage. Set the values first, then copy the number
# formatting from donor cells in row 14 that already carry the right style
# (C/D/F14 = text style, G14 = integer style, H14 = percent style) so the
# styles used match exactly.
foreach ($row in 15, 26) {
    $ws.Range("C$row").Value = "'0"
    $ws.Range("D$row").Value = 1
    $ws.Range("E$row").Value = -100

    $ws.Range("D14").Copy() | Out-Null
    $ws.Range("C$row").PasteSpecial(-4122) | Out-Null

    $ws.Range("G14").Copy() | Out-Null
    $ws.Range("D$row").PasteSpecial(-4122) | Out-Null

    $ws.Range("H14").Copy() | Out-Null
    $ws.Range("E$row").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# --- Remaining numeric value-only updates ---
$ws.Range("G14").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 14
$ws.Range("K15").Value = -28.571428571428
$ws.Range("N15").Value = -66.666666666666
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -46.153846153846
$ws.Range("I16").Value = 72
$ws.Range("J16").Value = 67
$ws.Range("K16").Value = 7.462686567164
$ws.Range("L16").Value = -2.702702702702
$ws.Range("M16").Value = -12.195121951219
$ws.Range("N16").Value = -84.175824175824
$ws.Range("C17").Value = 4
$ws.Range("E17").Value = -42.857142857142
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 9.523809523809
$ws.Range("I17").Value = 174
$ws.Range("J17").Value = 174
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = -15.121951219512
$ws.Range("M17").Value = 62.616822429906
$ws.Range("N17").Value = -52.717391304347
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 14.285714285714
$ws.Range("I18").Value = 52
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = 4
$ws.Range("L18").Value = 10.638297872340
$ws.Range("M18").Value = -46.938775510204
$ws.Range("N18").Value = -90.861159929701
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 400
$ws.Range("F19").Value = 11
$ws.Range("G19").Value = 15
$ws.Range("H19").Value = -26.666666666666
$ws.Range("I19").Value = 127
$ws.Range("J19").Value = 130
$ws.Range("K19").Value = -2.307692307692
$ws.Range("L19").Value = 12.389380530973
$ws.Range("M19").Value = 86.764705882352
$ws.Range("N19").Value = -44.052863436123
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -60
$ws.Range("I20").Value = 47
$ws.Range("J20").Value = 44
$ws.Range("K20").Value = 6.818181818181
$ws.Range("L20").Value = 46.875
$ws.Range("M20").Value = -24.193548387096
$ws.Range("N20").Value = -85.493827160493
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 54
$ws.Range("G21").Value = 70
$ws.Range("H21").Value = -22.857142857142
$ws.Range("I21").Value = 482
$ws.Range("J21").Value = 486
$ws.Range("K21").Value = -0.823045267489
$ws.Range("L21").Value = -2.626262626262
$ws.Range("M21").Value = 11.316397228637
$ws.Range("N21").Value = -75.730110775428
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 5
$ws.Range("I22").Value = 14
$ws.Range("K22").Value = 133.333333333333
$ws.Range("L22").Value = 133.333333333333
$ws.Range("M22").Value = 100
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 8
$ws.Range("H23").Value = 14.285714285714
$ws.Range("I23").Value = 65
$ws.Range("J23").Value = 61
$ws.Range("K23").Value = 6.557377049180
$ws.Range("L23").Value = -9.722222222222
$ws.Range("M23").Value = 150
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 16.666666666666
$ws.Range("F24").Value = 44
$ws.Range("G24").Value = 41
$ws.Range("H24").Value = 7.317073170731
$ws.Range("I24").Value = 398
$ws.Range("J24").Value = 415
$ws.Range("K24").Value = -4.096385542168
$ws.Range("L24").Value = -7.226107226107
$ws.Range("M24").Value = 49.624060150375
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 80
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = 12.121212121212
$ws.Range("I25").Value = 307
$ws.Range("J25").Value = 278
$ws.Range("K25").Value = 10.431654676259
$ws.Range("L25").Value = -0.647249190938
$ws.Range("M25").Value = 6.968641114982
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 25
$ws.Range("K26").Value = -12
$ws.Range("G27").Value = 4
$ws.Range("J27").Value = 36
$ws.Range("K27").Value = -19.444444444444
